$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 306; this pushes the existing
# rows 306-411 down to 307-412 (and grows the used range to row 412),
# matching the dimension change from A1:R411 to A1:R412.
$ws.Rows("306:306").Insert()

# Populate the newly inserted row with its data.
$ws.Range("A306").Value = 3
$ws.Range("B306").Value = "Femacal de La Calera"
$ws.Range("C306").Value = "Coquimbo"
$ws.Range("D306").Value = 44809
$ws.Range("E306").Value = 5
$ws.Range("F306").Value = 100114013
$ws.Range("G306").Value = "Zanahoria"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 310
$ws.Range("K306").Value = 11500
$ws.Range("L306").Value = 12000
$ws.Range("M306").Value = 11742
$ws.Range("N306").Value = "$/saco 20 kilos"
$ws.Range("O306").Value = "Provincia de Quillota"
$ws.Range("P306").Value = 587
$ws.Range("Q306").Value = 20
$ws.Range("R306").Value = "Hortaliza"
